# Update gh-pages to output generated at 456a3b4
# Refresh the "想去人数" (want-to-go count) column F on each sheet,
# plus one refreshed cover-image URL (I16 on 展览).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 312
$ws.Range("F7").Value = 5503
$ws.Range("F10").Value = 3782
$ws.Range("F12").Value = 18
$ws.Range("F13").Value = 21
$ws.Range("I16").Value = "//i1.hdslb.com/bfs/openplatform/202407/PRico1em1722234083280.jpeg"
$ws.Range("F18").Value = 98
$ws.Range("F20").Value = 167
$ws.Range("F21").Value = 126
$ws.Range("F23").Value = 5171
$ws.Range("F25").Value = 2058
$ws.Range("F27").Value = 329
$ws.Range("F28").Value = 7656
$ws.Range("F30").Value = 171
$ws.Range("F31").Value = 2178
$ws.Range("F32").Value = 2142
$ws.Range("F33").Value = 1321
$ws.Range("F35").Value = 1170
$ws.Range("F38").Value = 256
$ws.Range("F40").Value = 240
$ws.Range("F41").Value = 11
$ws.Range("F42").Value = 1173
$ws.Range("F44").Value = 27
$ws.Range("F45").Value = 1309
$ws.Range("F46").Value = 2011
$ws.Range("F47").Value = 112
$ws.Range("F48").Value = 205
$ws.Range("F49").Value = 1208

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 144
$ws.Range("F11").Value = 115
$ws.Range("F19").Value = 3

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 545
$ws.Range("F3").Value = 721
$ws.Range("F4").Value = 62

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 545
$ws.Range("F7").Value = 721
$ws.Range("F8").Value = 312
$ws.Range("F9").Value = 5503
$ws.Range("F10").Value = 3782
$ws.Range("F12").Value = 18
$ws.Range("F13").Value = 21
$ws.Range("F17").Value = 98
$ws.Range("F19").Value = 144
$ws.Range("F20").Value = 167
$ws.Range("F22").Value = 126
$ws.Range("F24").Value = 5171
$ws.Range("F26").Value = 2058
$ws.Range("F28").Value = 329
$ws.Range("F29").Value = 7657
$ws.Range("F31").Value = 171
$ws.Range("F32").Value = 2178
$ws.Range("F33").Value = 2142
$ws.Range("F34").Value = 1321
$ws.Range("F36").Value = 1170
$ws.Range("F37").Value = 256
$ws.Range("F38").Value = 240
$ws.Range("F39").Value = 11
$ws.Range("F40").Value = 1173
$ws.Range("F42").Value = 27
$ws.Range("F43").Value = 1309
$ws.Range("F45").Value = 2011
$ws.Range("F46").Value = 112
$ws.Range("F48").Value = 205
$ws.Range("F49").Value = 1208
